$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Proofread / rewrite the intro paragraph (2nd paragraph).
# ------------------------------------------------------------------
$oldIntro = "Login feature is required by client to provide security to users of the application. The users that are required to login are scorers, markers, judges and administrator ((if there is one). Each user will need a user name or user id and a password. Both must be validated and if it is valid checked against the database for a match. In case of a match users is successfully logged in. Otherwise an error message will be shown on the page accordingly."
$newIntro = "The login feature is required by the client to provide security to the users of the application. The users that are required to login are scorers, markers or the judge and the administrator (if there is one). Each user will need a user name/id and a password. Both must be validated and if it is valid checked against the database for a match. In case of a match users are successfully logged in. Otherwise an error message will be shown on the page accordingly."

$rng = $d.Content
$rng.Find.Execute($oldIntro, $true, $false, $false, $false, $false, $true, 1, $false, $newIntro, 2)

# ------------------------------------------------------------------
# 2. Move the (hidden) "_GoBack" bookmark from the top of the document
#    (just before the "Login Model" heading run) down to the end of
#    the last "Purpose" paragraph (right after the final run, before
#    the paragraph mark) -- this mirrors where Word itself drops
#    "_GoBack" after the most recent edit point.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Paragraph 5 is "The purpose of the login domain ... preserved and secured."
$p5 = $d.Paragraphs.Item(5)
$insertPos = $p5.Range.End - 1

# Placing a bookmark with Start == End exactly on the paragraph's last
# text boundary is unreliable, so stage a small unique marker there,
# wrap the bookmark around it, then collapse the marker away again --
# this leaves the bookmark correctly collapsed in the same spot.
$marker = "@@GOBACKMARK@@"
$insertionRange = $d.Range($insertPos, $insertPos)
$insertionRange.InsertAfter($marker)

$markerRange = $d.Range($insertPos, $insertPos + $marker.Length)
$d.Bookmarks.Add("_GoBack", $markerRange)

$clearRange = $d.Range($insertPos, $insertPos + $marker.Length)
$clearRange.Text = ""

Write-Output "Edit complete"
